# Update workbook to add carjacking data through 2022-03-24 (commit: "Add data for 2022-04-01")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename sheet and update the "March 2022" column header text ---
$ws.Name = "Through 2022-03-24"
$ws.Range("B1").Value = "March 2022 (through March 24)"

# --- Increment existing counts (new carjacking(s) recorded against these neighborhoods) ---
$ws.Range("B2").Value = 3    # Chinatown
$ws.Range("B3").Value = 9    # Austin
$ws.Range("Q5").Value = 9    # Garfield Park
$ws.Range("T9").Value = 2    # Chicago Lawn
$ws.Range("E14").Value = 2   # West Town
$ws.Range("H15").Value = 2   # Humboldt Park
$ws.Range("Q17").Value = 2   # Auburn Gresham
$ws.Range("E25").Value = 5   # Grand Crossing
$ws.Range("E26").Value = 3   # Chatham
$ws.Range("T36").Value = 3   # Roseland
$ws.Range("T49").Value = 2   # Little Village
$ws.Range("K77").Value = 3   # Portage Park

# --- Remove a stale data point (re-categorized away from West Loop / March 2022) ---
$ws.Range("B10").ClearContents()   # West Loop

# --- Add newly-populated cells (value 1) ---
$ws.Range("Q13").Value = 1   # Woodlawn
$ws.Range("K19").Value = 1   # Lincoln Park
$ws.Range("H24").Value = 1   # Washington Park
$ws.Range("Q26").Value = 1   # Chatham
$ws.Range("B28").Value = 1   # Logan Square
$ws.Range("T29").Value = 1   # Lower West Side
$ws.Range("K30").Value = 1   # United Center
$ws.Range("E34").Value = 1   # River North
$ws.Range("Q88").Value = 1   # Uptown

$wb.Save()
